$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38 (shifts rows 38..48 down to 39..49)
$ws.Rows.Item(38).Insert()

# Fill in the new row's content, matching the style of the surrounding normal rows (s="2")
$ws.Range("A38").Value = "为文章指定封面图片"
$ws.Range("A38").Style = $ws.Range("A40").Style

# Update the selection / view to match the committed file
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D35").Select()
